$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update purpose (column E, rows 2-54) from "Retrofitted_2496" to "fullRNASEQ"
$ws.Range("E2:E54").Value = "fullRNASEQ"

# Update libraryPreparer (column B, rows 2-54) from "Retrofitted_2496" to "H.BROWN"
$ws.Range("B2:B54").Value = "H.BROWN"

# Reflect the selection left behind in the sheet after editing column B
$ws.Range("B3:B54").Select()
